# Draft of the templates for the function for reloading historical data
# during FHIT import.
#
# Inserts a new row (44) above the existing "25_..." template block so a
# new template row describing 25_adding_historical_raw_records.sql /
# template_adding_historical_records.sql is added, pushing all the rows
# below it (old 44-59) down by one (new 45-60). Also relocates the two
# cell comments that lived in the shifted area and restores the
# selection/scroll position shown in the author's Excel session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Preserve the two comments that sit in the area which is about to be
#    shifted down, before we touch any rows.
# ---------------------------------------------------------------------
$commentText1 = $ws.Range("K50").Comment.Text()
$commentText2 = $ws.Range("K52").Comment.Text()
$ws.Range("K50").Comment.Delete()
$ws.Range("K52").Comment.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new row at position 44 - this shifts existing rows 44-59
#    down to 45-60 (values, formulas and styles move with the cells).
# ---------------------------------------------------------------------
$ws.Rows("44:44").Insert()

# ---------------------------------------------------------------------
# 3. Populate the newly inserted row 44 with the new template entry.
# ---------------------------------------------------------------------
$ws.Range("B44").Value = "25_adding_historical_raw_records.sql"
$ws.Range("C44").Value = "template_adding_historical_records.sql"
$ws.Range("D44").Value = "cds2db_user"
$ws.Range("E44").Value = "cds2db_in"
$ws.Range("G44").Value = "v_"
$ws.Range("H44").Value = "_raw_last_version"
$ws.Range("J44").Value = "cds2db_user"
$ws.Range("N44").Value = "cds2db_in"
$ws.Range("O44").Value = "_raw"

# J44 keeps the same "plain" override style used by the analogous cells
# above it (J33/J38/J39/J40/J41) instead of the default inherited look.
$ws.Range("J41").Copy()
$ws.Range("J44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Re-create the two comments at their new (shifted) locations.
# ---------------------------------------------------------------------
$ws.Range("K51").AddComment($commentText1)
$ws.Range("K53").AddComment($commentText2)

# ---------------------------------------------------------------------
# 5. Restore view state: scrolled position and active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("J44").Select()
